$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.601781964302063
$ws.Range("B1").Value = 1.670974373817444
$ws.Range("C1").Value = 1.74342405796051
$ws.Range("D1").Value = 2.353389978408813
$ws.Range("E1").Value = 3.929531335830688
